# Milestone4 Presentation3 - "Next Steps (1)" slide:
#  - insert a new bullet "Improve the rotation function" between the
#    "...java function" sub-bullet and "getting the ?task(...)" bullet
#  - set explicit font sizes on every paragraph (24pt for top-level
#    bullets, 22pt for the 2nd level, 21pt for the 3rd level) now that
#    the extra line no longer lets autofit pick the sizes implicitly

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# The new bullet belongs right before "getting the ?task(_,_,_,_) ..."
# (currently paragraph 3) and should sit at the top outline level, same
# as that paragraph - so insert before it and the level is inherited.
$thirdPara = $tr.Paragraphs(3)
[void]$thirdPara.InsertBefore("Improve the rotation function`r")

# Re-fetch after the text changed and apply the per-level font sizes.
$tr = $sh.TextFrame.TextRange
$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i)
    $lvl = $para.IndentLevel
    if ($lvl -eq 1) {
        $para.Font.Size = 24
    } elseif ($lvl -eq 2) {
        $para.Font.Size = 22
    } elseif ($lvl -eq 3) {
        $para.Font.Size = 21
    }
}
